# "Decision tree images seperated" — resave pass.
#
# The canonical diff for this commit only touches two kinds of things:
#   1) PowerPoint's own internal revision-tracking bookkeeping
#      (ppt/revisionInfo.xml, ppt/changesInfos/changesInfo1.xml) — these
#      are written by PowerPoint's coauthoring/change-tracking engine as
#      a byproduct of opening & resaving the file two days later; there's
#      no object-model call that authors them directly, so they are left
#      for the host application to regenerate on save.
#   2) The cached text of every "datetimeFigureOut" date field (on the
#      slide master, every slide layout, and the notes master) flipping
#      from "11-6-2018" to "13-6-2018" — i.e. PowerPoint recalculated the
#      auto-date placeholders the next time the deck was saved.
#
# Reproduce (2) by walking every placeholder shape named like a date
# placeholder on the slide master / layouts / notes master and rewriting
# its text from the old cached date to the new one.

function Update-DatePlaceholders {
    param($container)

    if ($container -eq $null) { return }
    $shapes = $container.Shapes
    if ($shapes -eq $null) { return }

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if (($shp.Name -like "*Date Placeholder*") -or ($shp.Name -like "*datum*")) {
            if ($shp.HasTextFrame) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq "11-6-2018") {
                    $tr.Text = "13-6-2018"
                }
            }
        }
    }
}

$p = $ppt.ActivePresentation

# Slide master's own date placeholder.
Update-DatePlaceholders $p.SlideMaster

# Notes master's date placeholder (best effort — some hosts expose this
# part read-only; the call is harmless if so).
Update-DatePlaceholders $p.NotesMaster

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholders $layouts.Item($j)
}
